{"js": "// Auto update 2026-02-09 15:05:09\n// Applies the dissertation-proposal text edits: refreshed timestamp,\n// wording tweaks, a citation-year fix, and reference-list DOI/URL additions\n// (plus one brand-new reference paragraph at the end of the list).\nconst body = context.document.body;\n\n// [oldText, newText] pairs - each oldText is unique & verbatim in the doc,\n// except the \"(Roth et al., ...; Defard et al., 2020)\" pair which occurs twice\n// (Background section + 4.2 Anomaly Detection Baselines) and must be fixed both times.\nconst replacements = [\n  [\"Generated: 2026-02-09 15:03\", \"Generated: 2026-02-09 15:05\"],\n  [\"This proposal outlines the design and implementation of an industrial image anomaly inspection pipeline with interpretable semantics and safety-aligned risk scoring. The system targets MVTec AD categories (bottle, cable, wood, tile, leather) and outputs anomaly detection, localization, defect labels, evidence, risk class, and action recommendations. The work combines state-of-the-art anomaly detection baselines with a vision-language model for defect semantics and deterministic risk rules.\", \"This proposal outlines the design and implementation of an industrial image anomaly inspection pipeline with interpretable semantics and safety-aligned risk scoring. The system targets MVTec AD categories (bottle, cable, wood, tile, leather) and outputs anomaly detection, localization, defect labels, evidence, risk class, and action recommendations. The work combines established anomaly detection baselines with a vision-language model for defect semantics and deterministic risk rules.\"],\n  [\"Industrial visual inspection often relies on scarce defect examples, making unsupervised or one-class anomaly detection appropriate. The MVTec AD dataset provides a realistic benchmark with pixel-precise annotations for industrial anomalies (Bergmann et al., 2019; Bergmann et al., 2021). Recent methods such as PatchCore and PaDiM demonstrate strong performance for detection and localization in this setting (Roth et al., 2021; Defard et al., 2020).\", \"Industrial visual inspection often relies on scarce defect examples, making unsupervised or one-class anomaly detection appropriate. The MVTec AD dataset provides a realistic benchmark with pixel-precise annotations for industrial anomalies (Bergmann et al., 2019; Bergmann et al., 2021). Recent methods such as PatchCore and PaDiM demonstrate strong performance for detection and localization in this setting (Roth et al., 2022; Defard et al., 2020).\"],\n  [\"Implement PatchCore and PaDiM as strong baselines for detection and localization. PatchCore uses a memory bank of nominal patch features for outlier scoring, while PaDiM models patch embeddings with multivariate Gaussians (Roth et al., 2021; Defard et al., 2020).\", \"Implement PatchCore and PaDiM as strong baselines for detection and localization. PatchCore uses a memory bank of nominal patch features for outlier scoring, while PaDiM models patch embeddings with multivariate Gaussians (Roth et al., 2022; Defard et al., 2020).\"],\n  [\"Integrate LLaVA for constrained defect labeling. LLaVA is a multimodal model created via visual instruction tuning, providing strong image-language reasoning (Liu et al., 2023). The LLaVA-1.6 release improves visual reasoning and OCR, supporting detailed defect evidence when constrained to fixed label sets (Liu et al., 2024).\", \"Integrate LLaVA for constrained defect labeling. LLaVA is a multimodal model created via visual instruction tuning, providing strong image-language reasoning (Liu et al., 2023). The LLaVA-1.6 (Mistral) release improves visual reasoning and supports higher-resolution inputs, which is useful for defect evidence extraction when constrained to fixed label sets (Liu et al., 2024).\"],\n  [\"Bergmann, P., Fauser, M., Sattlegger, D. and Steger, C. (2019) MVTec AD \u2014 A Comprehensive Real-World Dataset for Unsupervised Anomaly Detection. Proceedings of the IEEE/CVF Conference on Computer Vision and Pattern Recognition.\", \"Bergmann, P., Fauser, M., Sattlegger, D. and Steger, C. (2019) MVTec AD \u2014 A Comprehensive Real-World Dataset for Unsupervised Anomaly Detection. Proceedings of the IEEE/CVF Conference on Computer Vision and Pattern Recognition. doi:10.1109/CVPR.2019.00982.\"],\n  [\"Bergmann, P., Batzner, K., Fauser, M., Sattlegger, D. and Steger, C. (2021) The MVTec Anomaly Detection Dataset: A Comprehensive Real-World Dataset for Unsupervised Anomaly Detection. International Journal of Computer Vision.\", \"Bergmann, P., Batzner, K., Fauser, M., Sattlegger, D. and Steger, C. (2021) The MVTec Anomaly Detection Dataset: A Comprehensive Real-World Dataset for Unsupervised Anomaly Detection. International Journal of Computer Vision, 129(4), pp.1038\u20131059. doi:10.1007/s11263-020-01400-4.\"],\n  [\"Defard, T., Setkov, A., Loesch, A. and Audigier, R. (2020) PaDiM: a Patch Distribution Modeling Framework for Anomaly Detection and Localization. arXiv preprint arXiv:2011.08785.\", \"Defard, T., Setkov, A., Loesch, A. and Audigier, R. (2020) PaDiM: a Patch Distribution Modeling Framework for Anomaly Detection and Localization. arXiv:2011.08785. doi:10.48550/arXiv.2011.08785.\"],\n  [\"Roth, K., Pemula, L., Zepeda, J., Sch\u00f6lkopf, B., Brox, T. and Gehler, P. (2021) Towards Total Recall in Industrial Anomaly Detection. arXiv preprint arXiv:2106.08265.\", \"Roth, K., Pemula, L., Zepeda, J., Sch\u00f6lkopf, B., Brox, T. and Gehler, P. (2022) Towards Total Recall in Industrial Anomaly Detection. arXiv:2106.08265. doi:10.48550/arXiv.2106.08265.\"],\n  [\"Liu, H., Li, C., Wu, Q. and Lee, Y.J. (2023) Visual Instruction Tuning. arXiv preprint arXiv:2304.08485.\", \"Liu, H., Li, C., Wu, Q. and Lee, Y.J. (2023) Visual Instruction Tuning. arXiv:2304.08485. doi:10.48550/arXiv.2304.08485.\"],\n  [\"Liu, H., Li, C., Li, Y., Li, B., Zhang, Y., Shen, S. and Lee, Y.J. (2024) LLaVA-1.6: Improved reasoning, OCR, and world knowledge (blog release).\", \"Liu, H., Li, C., Li, Y., Li, B., Zhang, Y., Shen, S. and Lee, Y.J. (2024) LLaVA-1.6 (Mistral) model card. Hugging Face. Available at: https://huggingface.co/llava-hf/llava-v1.6-mistral-7b-hf (Accessed: 9 February 2026).\"],\n];\n\nconst newParagraphText = \"MVTec Software (n.d.) MVTec AD dataset page. Available at: https://www.mvtec.com/company/research/datasets/mvtec-ad (Accessed: 9 February 2026).\";\nconst lastRefText = \"Liu, H., Li, C., Li, Y., Li, B., Zhang, Y., Shen, S. and Lee, Y.J. (2024) LLaVA-1.6 (Mistral) model card. Hugging Face. Available at: https://huggingface.co/llava-hf/llava-v1.6-mistral-7b-hf (Accessed: 9 February 2026).\";\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found for replacement: \" + oldText.slice(0, 60));\n  }\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// Append the new reference paragraph after the (now updated) LLaVA-1.6 reference,\n// i.e. at the very end of the References list, matching the style of its neighbors.\nconst lastRefResults = body.search(lastRefText, { matchCase: true });\nlastRefResults.load(\"items\");\nawait context.sync();\nif (lastRefResults.items.length === 0) {\n  throw new Error(\"Could not locate the final reference paragraph to insert after.\");\n}\nconst lastRefRange = lastRefResults.items[0];\nconst lastRefParagraph = lastRefRange.paragraphs.getFirst();\nlastRefParagraph.insertParagraph(newParagraphText, \"After\");\nawait context.sync();\n\n", "ps1": "# Auto update 2026-02-09 15:05:09\n# Applies the dissertation-proposal text edits: refreshed timestamp,\n# wording tweaks, a citation-year fix, and reference-list DOI/URL additions\n# (plus one brand-new reference paragraph at the end of the list).\n$d = $word.ActiveDocument\n\n# [oldText, newText] pairs - each oldText is unique & verbatim in the doc,\n# except the \"(Roth et al., ...; Defard et al., 2020)\" pair which occurs twice\n# (Background section + 4.2 Anomaly Detection Baselines); Find/Replace (wdReplaceAll)\n# fixes both occurrences in one call.\n$replacements = @(\n    @('Generated: 2026-02-09 15:03', 'Generated: 2026-02-09 15:05'),\n    @('This proposal outlines the design and implementation of an industrial image anomaly inspection pipeline with interpretable semantics and safety-aligned risk scoring. The system targets MVTec AD categories (bottle, cable, wood, tile, leather) and outputs anomaly detection, localization, defect labels, evidence, risk class, and action recommendations. The work combines state-of-the-art anomaly detection baselines with a vision-language model for defect semantics and deterministic risk rules.', 'This proposal outlines the design and implementation of an industrial image anomaly inspection pipeline with interpretable semantics and safety-aligned risk scoring. The system targets MVTec AD categories (bottle, cable, wood, tile, leather) and outputs anomaly detection, localization, defect labels, evidence, risk class, and action recommendations. The work combines established anomaly detection baselines with a vision-language model for defect semantics and deterministic risk rules.'),\n    @('Industrial visual inspection often relies on scarce defect examples, making unsupervised or one-class anomaly detection appropriate. The MVTec AD dataset provides a realistic benchmark with pixel-precise annotations for industrial anomalies (Bergmann et al., 2019; Bergmann et al., 2021). Recent methods such as PatchCore and PaDiM demonstrate strong performance for detection and localization in this setting (Roth et al., 2021; Defard et al., 2020).', 'Industrial visual inspection often relies on scarce defect examples, making unsupervised or one-class anomaly detection appropriate. The MVTec AD dataset provides a realistic benchmark with pixel-precise annotations for industrial anomalies (Bergmann et al., 2019; Bergmann et al., 2021). Recent methods such as PatchCore and PaDiM demonstrate strong performance for detection and localization in this setting (Roth et al., 2022; Defard et al., 2020).'),\n    @('Implement PatchCore and PaDiM as strong baselines for detection and localization. PatchCore uses a memory bank of nominal patch features for outlier scoring, while PaDiM models patch embeddings with multivariate Gaussians (Roth et al., 2021; Defard et al., 2020).', 'Implement PatchCore and PaDiM as strong baselines for detection and localization. PatchCore uses a memory bank of nominal patch features for outlier scoring, while PaDiM models patch embeddings with multivariate Gaussians (Roth et al., 2022; Defard et al., 2020).'),\n    @('Integrate LLaVA for constrained defect labeling. LLaVA is a multimodal model created via visual instruction tuning, providing strong image-language reasoning (Liu et al., 2023). The LLaVA-1.6 release improves visual reasoning and OCR, supporting detailed defect evidence when constrained to fixed label sets (Liu et al., 2024).', 'Integrate LLaVA for constrained defect labeling. LLaVA is a multimodal model created via visual instruction tuning, providing strong image-language reasoning (Liu et al., 2023). The LLaVA-1.6 (Mistral) release improves visual reasoning and supports higher-resolution inputs, which is useful for defect evidence extraction when constrained to fixed label sets (Liu et al., 2024).'),\n    @('Bergmann, P., Fauser, M., Sattlegger, D. and Steger, C. (2019) MVTec AD \u2014 A Comprehensive Real-World Dataset for Unsupervised Anomaly Detection. Proceedings of the IEEE/CVF Conference on Computer Vision and Pattern Recognition.', 'Bergmann, P., Fauser, M., Sattlegger, D. and Steger, C. (2019) MVTec AD \u2014 A Comprehensive Real-World Dataset for Unsupervised Anomaly Detection. Proceedings of the IEEE/CVF Conference on Computer Vision and Pattern Recognition. doi:10.1109/CVPR.2019.00982.'),\n    @('Bergmann, P., Batzner, K., Fauser, M., Sattlegger, D. and Steger, C. (2021) The MVTec Anomaly Detection Dataset: A Comprehensive Real-World Dataset for Unsupervised Anomaly Detection. International Journal of Computer Vision.', 'Bergmann, P., Batzner, K., Fauser, M., Sattlegger, D. and Steger, C. (2021) The MVTec Anomaly Detection Dataset: A Comprehensive Real-World Dataset for Unsupervised Anomaly Detection. International Journal of Computer Vision, 129(4), pp.1038\u20131059. doi:10.1007/s11263-020-01400-4.'),\n    @('Defard, T., Setkov, A., Loesch, A. and Audigier, R. (2020) PaDiM: a Patch Distribution Modeling Framework for Anomaly Detection and Localization. arXiv preprint arXiv:2011.08785.', 'Defard, T., Setkov, A., Loesch, A. and Audigier, R. (2020) PaDiM: a Patch Distribution Modeling Framework for Anomaly Detection and Localization. arXiv:2011.08785. doi:10.48550/arXiv.2011.08785.'),\n    @('Roth, K., Pemula, L., Zepeda, J., Sch\u00f6lkopf, B., Brox, T. and Gehler, P. (2021) Towards Total Recall in Industrial Anomaly Detection. arXiv preprint arXiv:2106.08265.', 'Roth, K., Pemula, L., Zepeda, J., Sch\u00f6lkopf, B., Brox, T. and Gehler, P. (2022) Towards Total Recall in Industrial Anomaly Detection. arXiv:2106.08265. doi:10.48550/arXiv.2106.08265.'),\n    @('Liu, H., Li, C., Wu, Q. and Lee, Y.J. (2023) Visual Instruction Tuning. arXiv preprint arXiv:2304.08485.', 'Liu, H., Li, C., Wu, Q. and Lee, Y.J. (2023) Visual Instruction Tuning. arXiv:2304.08485. doi:10.48550/arXiv.2304.08485.'),\n    @('Liu, H., Li, C., Li, Y., Li, B., Zhang, Y., Shen, S. and Lee, Y.J. (2024) LLaVA-1.6: Improved reasoning, OCR, and world knowledge (blog release).', 'Liu, H., Li, C., Li, Y., Li, B., Zhang, Y., Shen, S. and Lee, Y.J. (2024) LLaVA-1.6 (Mistral) model card. Hugging Face. Available at: https://huggingface.co/llava-hf/llava-v1.6-mistral-7b-hf (Accessed: 9 February 2026).'),\n)\n\n$newParagraphText = 'MVTec Software (n.d.) MVTec AD dataset page. Available at: https://www.mvtec.com/company/research/datasets/mvtec-ad (Accessed: 9 February 2026).'\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # wdReplaceAll (2) so the duplicated Roth/Defard citation is fixed everywhere it occurs.\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found for replacement: $oldText\"\n    }\n}\n\n# Append the new reference paragraph at the very end of the References list,\n# right after the (now updated) LLaVA-1.6 reference.\n$lastParaIndex = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($lastParaIndex)\n$lastPara.Range.InsertParagraphAfter() | Out-Null\n\n$newLastParaIndex = $d.Paragraphs.Count\n$newPara = $d.Paragraphs.Item($newLastParaIndex)\n$newPara.Range.Text = $newParagraphText\n\n"}
